# Replace "yomeru"/"yomenai" (read) with "mieru"/"mienai" (see/look) text,
# and fix "horror-de-nai eiga ga mitai" to use the kanji for "mitai".
#
# Three shapes on slide 1 are touched:
#   id=47 "テキスト ボックス 46" : "：読める"   -> "：" + "見える"   (split into two runs)
#   id=49 "テキスト ボックス 48" : "：読めない" -> "：" + "見えない" (split into two runs)
#   id=56 "テキスト ボックス 55" : "ホラーでない映画がみたい" -> "ホラーでない映画が見たい"
#       (2nd paragraph, plain text substitution - stays a single run)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# --- Shape 47: "：読める" -> "：見える" (split run after the full-width colon) ---
$shp47 = Get-ShapeById $s 47
$tr47 = $shp47.TextFrame.TextRange
# "：読める" -> chars: 1=："，2=読，3=め，4=る ; replace "読める" (2,3) with "見える"
$run47 = $tr47.Characters(2, 3)
$run47.Text = "見える"

# --- Shape 49: "：読めない" -> "：見えない" (split run after the full-width colon) ---
$shp49 = Get-ShapeById $s 49
$tr49 = $shp49.TextFrame.TextRange
# "：読めない" -> chars: 1=：，2=読，3=め，4=な，5=い ; replace "読めない" (2,4) with "見えない"
$run49 = $tr49.Characters(2, 4)
$run49.Text = "見えない"

# --- Shape 56: 2nd paragraph "ホラーでない映画がみたい" -> "...見たい" (no run split) ---
$shp56 = Get-ShapeById $s 56
$tr56 = $shp56.TextFrame.TextRange
# Full text is "SFかアニメか、" + vertical-tab paragraph break + "ホラーでない映画がみたい"
# (8 chars) + (1 break char) = the 2nd paragraph's run starts at char 10, length 12.
$run56 = $tr56.Characters(10, 12)
$run56.Text = "ホラーでない映画が見たい"
